$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 582 (pushes existing rows 582-646 down to 585-649)
$ws.Rows.Item(582).Resize(3).Insert()

# New weekly data block (date 2022-12-23 -> serial 44918)
$newDate = 44918

$rows = @(582, 583, 584)
$calidad = @("Maduro", "Pintón", "Primera Pintón")
$volumen = @(200, 400, 400)
$precioMin = @(18000, 19000, 21000)
$precioMax = @(18000, 19000, 21000)
$precioProm = @(18000, 19000, 21000)
$precioKg = @(900, 950, 1050)

for ($i = 0; $i -lt 3; $i++) {
    $r = $rows[$i]

    $ws.Cells.Item($r, 1).Value = 11
    $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($r, 3).Value = "Bíobío"
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 8
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108006
    $ws.Cells.Item($r, 10).Value = "Plátano"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $calidad[$i]
    $ws.Cells.Item($r, 13).Value = $volumen[$i]
    $ws.Cells.Item($r, 14).Value = $precioMin[$i]
    $ws.Cells.Item($r, 15).Value = $precioMax[$i]
    $ws.Cells.Item($r, 16).Value = $precioProm[$i]
    $ws.Cells.Item($r, 17).Value = "`$/caja 20 kilos"
    $ws.Cells.Item($r, 18).Value = "Ecuador"
    $ws.Cells.Item($r, 19).Value = $precioKg[$i]
    $ws.Cells.Item($r, 20).Value = 20
}
